$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2 through 42: Column A = Actual Consumption (MW), Column B = Timestamp (serial date)
$data = @(
    @(5047, 45733),
    @(5011, 45733.01041666666),
    @(4991, 45733.02083333334),
    @(4939, 45733.03125),
    @(4914, 45733.04166666666),
    @(4892, 45733.05208333334),
    @(4892, 45733.0625),
    @(4884, 45733.07291666666),
    @(4937, 45733.08333333334),
    @(4939, 45733.09375),
    @(4922, 45733.10416666666),
    @(4975, 45733.11458333334),
    @(5044, 45733.125),
    @(5015, 45733.13541666666),
    @(5038, 45733.14583333334),
    @(5139, 45733.15625),
    @(5235, 45733.16666666666),
    @(5214, 45733.17708333334),
    @(5359, 45733.1875),
    @(5519, 45733.19791666666),
    @(5750, 45733.20833333334),
    @(5824, 45733.21875),
    @(6006, 45733.22916666666),
    @(6215, 45733.23958333334),
    @(6517, 45733.25),
    @(6680, 45733.26041666666),
    @(6867, 45733.27083333334),
    @(7012, 45733.28125),
    @(7080, 45733.29166666666),
    @(7221, 45733.30208333334),
    @(7314, 45733.3125),
    @(7346, 45733.32291666666),
    @(7268, 45733.33333333334),
    @(7337, 45733.34375),
    @(7359, 45733.35416666666),
    @(7222, 45733.36458333334),
    @(7179, 45733.375),
    @(7178, 45733.38541666666),
    @(7154, 45733.39583333334),
    @(7073, 45733.40625),
    @(7055, 45733.41666666666)
)

$dateFormat = $ws.Range("B40").NumberFormat

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 2).NumberFormat = $dateFormat
    $r++
}

Write-Host "Updated rows 2-42 with refreshed Entsoe consumption data"
